# Build site at 2023-04-12 14:53:07 UTC
# Applies the LOT2025 discipline-sheet update:
#  - inserts a new row (old row 13 "Programa resumido" block shifts down,
#    and the professor name gets its own dedicated row)
#  - fills in the new "Objetivos" text (PT)
#  - fills the professor name into its new row
#  - fills the new short/full syllabus text (PT)
#  - fills the Metodo / Criterio / Norma de recuperacao texts in their
#    (now correctly aligned) rows
#  - fills in the new bibliography text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 13 - shifts existing rows 13:24 down to 14:25
$ws.Rows("13:13").Insert()

# The insert carries column A's bold style into the new row; that row has no
# label in column A (only B/C get content), so clear it back out.
$ws.Range("A13").Clear()

# Row 10 (Objetivos:) - new Portuguese objectives text
$objetivos = "Capacitar o aluno para identificar e aplicar os conceitos de Reatores Bioquímicos em diferentes bioprocessos (fermentativos e enzimáticos). Especificamente, capacitar o aluno para definir os tipos de biorreatores; para definir as diferentes formas de condução de um processo fermentativo empregando biorreatores, para realizar o equacionamento matemático do crescimento microbiano e da formação de produtos de interesse de um bioprocesso empregando biorreatores e para aplicação dos conceitos gerais dos reatores enzimáticos operados em diferentes fases."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Row 13 (new row, under "Docentes responsaveis:") - professor name
$professor = "1112574 - Inês Conceição Roberto"
$ws.Range("B13").Value = $professor
$ws.Range("C13").Value = $professor
# Re-apply the normal (non-bold, wrapped) look that every other B-column cell uses
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160

# Row 14 (Programa resumido:) - new short syllabus text (PT)
$programaResumido = "Introdução a biorreatores; processo descontínuo; processo contínuo; processo descontínuo alimentado e reatores enzimáticos."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 16 (Programa:) - new full syllabus text (PT)
$programa = "1. Introdução a biorreatores: apresentação e classificação de reatores bioquímicos; interação microorganismos/meios para estabelecimento de condições para cálculos de biorreatores.2. Processo descontínuo: características gerais do processo contínuo; balanço material para célula, substrato e produto, em um único estágio com e sem reciclo de células; aplicação do processo contínuo (exemplos).3. Processo contínuo: características gerais do processo contínuo; balanço material para célula, substrato e produto, em um único estágio com e sem reciclo de células; aplicação do processo contínuo (exemplos).4. Processo descontínuo alimentado: características gerais do processo descontínuo alimentado; balanço material para célula e substrato, com volume variável, empregando vazão constante de alimentação; considerações sobre formação de produtos no processo descontínuo alimentado; aplicação do processo descontínuo alimentado (exemplos).5. Reatores enzimáticos: características gerais dos reatores enzimáticos; aplicação de processos enzimáticos (exemplos)."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Row 19 (Metodo:) - evaluation method text
$metodo = "Os alunos serão avaliados formalmente por duas provas teóricas. A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1 + P2 )/2"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 (Criterio:) - passing criterion text
$criterio = "Serão aprovados os alunos que obtiverem média igual ou maior que 5,0."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 (Norma de recuperacao:) - recovery norm text
$normaRecuperacao = "Aos alunos que não obtiverem média igual ou maior que 5,0, será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2 Serão aprovados os alunos que obtiverem média igual ou maior que 5,0"
$ws.Range("B21").Value = $normaRecuperacao
$ws.Range("C21").Value = $normaRecuperacao

# Row 22 (Bibliografia:) - new bibliography text
$bibliografia = "1)  Aiba, S., Humphrey, A.E., Millis, N.F. Biochemical Engineering - 2ª Edição- 1973.2) Asenjo A., Merchuk, J.C. Bioreactor System Design-1995.3) Stanbury, D. and Whitaker, A. Principles af Fermentation Technology-1986.4) Lima, U.A., Aquarone, E., Borzani, W. Biotecnologia Industrial. Fundamentos Vol. 1, Engenharia Bioquímica Vol.2, Processos Fermentativos Vol.3. Ed.  Edgard Blucher, São Paulo, 2001."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
